# V. 108 "Code 8"
# Add the movie "Code 8 (2019)" to the "Películas" ratings table (Tabla24).
# The table is sorted descending by column C ("Puntuación total"), so the
# new entry (average score 5.8) lands right after the current row 98
# (score 5.8125) and pushes every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")
$tbl = $ws.ListObjects.Item("Tabla24")

# Insert a fresh row at the correct sorted position (row 99); this shifts
# the existing rows 99:129 down to 100:130, carrying their values, formulas
# (auto-adjusted) and formatting with them.
$ws.Rows("99:99").Insert(-4121)

# Fill in the new movie's data.
$ws.Range("B99").Value2 = "Code 8 (2019)"
$ws.Range("C99").Formula = "=AVERAGE(D99,E99,E99,F99,G99,H99,H99,I99)"
$ws.Range("D99").Value2 = 7
$ws.Range("E99").Value2 = 6
$ws.Range("F99").Value2 = 5
$ws.Range("G99").Value2 = 5
$ws.Range("H99").Value2 = 6.1
$ws.Range("I99").Value2 = 5.2

# The worksheet-level row insert doesn't auto-grow the table, so extend the
# ListObject to cover the new row.
$tbl.Resize($ws.Range("B2:I130"))

# Re-point the visible selection at the new last row, as in the saved file.
$ws.Activate() | Out-Null
$ws.Range("C130").Select() | Out-Null
